$d = $word.ActiveDocument

# --- Step 1: remove the square brackets around every placeholder
# ([Customer], [Product] x5, [YourName]) while leaving the inner
# text / formatting / proofErr markers untouched.
$d.Content.Find.Execute("[", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
$d.Content.Find.Execute("]", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Step 2: paragraph "2. The Product operates effectively under
# various conditions." used to live in a single run (the brackets
# were removed above, but "Product" still needs to become bold/red
# like the other two "Product" mentions, and the editing cursor
# bookmark "_GoBack" needs to sit around it).
$p8 = $d.Paragraphs(8).Range
$pStart = $p8.Start

# wipe the paragraph's text, keeping its paragraph mark / pPr
$full = $d.Range($pStart, $p8.End - 1)
$full.Text = ""

$r1 = $d.Range($pStart, $pStart)
$r1.InsertAfter("2. The ")
$r1.LanguageID = "en-US"

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("Product")
$r2.Font.Bold = 1
$r2.Font.Color = 255
$r2.LanguageID = "en-US"

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(" ")
$r3.LanguageID = "en-US"

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter("operates effectively under various conditions.")
$r4.LanguageID = "en-US"

# moves (recreates) the "_GoBack" bookmark here, removing it from
# its old location automatically (bookmark names are unique)
$d.Bookmarks.Add("_GoBack", $r2)
